# Refresh the cryptos list with updated prices / 1h-volume percentages,
# and re-rank two coin pairs that swapped positions (Bittensor/Fetch.AI
# and Maker/Stacks), as published by the GitHub Actions refresh job.
#
# Values are written through a NumberFormat="@" / Style="Normal" round
# trip so price strings that look numeric (e.g. "1.00", "0.0000310",
# multi-dot figures like "70.302.65") are preserved verbatim as text
# instead of being auto-coerced into numbers (which would silently drop
# trailing zeros or flip into scientific notation) -- while leaving the
# cell's style back at the workbook's default ("Normal"/no explicit
# style index), matching the original formatting.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($CellRef, $Value) {
    $r = $ws.Range($CellRef)
    $r.NumberFormat = "@"
    $r.Value = $Value
    $r.Style = "Normal"
}

$updates = [ordered]@{
    "D2" = '70.302.65'
    "E2" = '  -1.53%  '
    "D3" = '3.549.01'
    "E3" = '  -1.43%  '
    "E4" = '  -0.09%  '
    "D5" = '618.38'
    "E5" = '  +5.32%  '
    "D6" = '188.02'
    "E6" = '  +0.27%  '
    "E7" = '  +2.42%  '
    "D8" = '1.00'
    "E8" = '  -0.03%  '
    "D9" = '0.217'
    "E9" = '  -1.37%  '
    "E10" = '  +1.26%  '
    "D11" = '53.93'
    "E11" = '  -1.26%  '
    "D12" = '0.0000310'
    "E12" = '  -4.53%  '
    "E13" = '  +1.46%  '
    "D14" = '4.116.14'
    "E14" = '  -1.22%  '
    "D15" = '617.74'
    "E15" = '  +8.16%  '
    "D16" = '70.358.76'
    "E16" = '  -1.30%  '
    "D17" = '12.82'
    "E17" = '  +3.25%  '
    "D18" = '19.15'
    "E18" = '  -1.22%  '
    "D19" = '3.535.87'
    "E19" = '  -1.26%  '
    "E20" = '  +0.05%  '
    "E21" = '  -1.46%  '
    "D22" = '17.73'
    "E22" = '  +0.46%  '
    "D23" = '105.64'
    "E23" = '  +10.92%  '
    "E24" = '  +2.32%  '
    "D25" = '5.13'
    "E25" = '  +0.25%  '
    "D26" = '3.04'
    "E26" = '  +3.31%  '
    "E27" = '  -3.80%  '
    "D28" = '10.05'
    "E28" = '  +9.37%  '
    "D29" = '34.59'
    "E29" = '  +6.60%  '
    "E30" = '  -2.42%  '
    "D31" = '12.56'
    "E31" = '  +1.68%  '
    "E32" = '  +1.54%  '
    "D33" = '64.21'
    "E33" = '  -0.27%  '
    "D34" = '3.76'
    "E34" = '  +15.07%  '
    "B35" = 'Fetch.AI'
    "C35" = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
    "D35" = '3.17'
    "E35" = '  -5.78%  '
    "B36" = 'Bittensor'
    "C36" = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
    "D36" = '537.59'
    "E36" = '  -3.20%  '
    "E37" = '  -0.06%  '
    "E38" = '  -3.76%  '
    "D39" = '37.41'
    "E39" = '  -1.08%  '
    "D40" = '0.0₃0784'
    "E40" = '  -3.77%  '
    "B41" = 'Stacks'
    "C41" = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
    "D41" = '3.57'
    "E41" = '  +3.60%  '
    "B42" = 'Maker'
    "C42" = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
    "D42" = '3.547.59'
    "E42" = '  +0.92%  '
    "E43" = '  +2.04%  '
    "D44" = '0.0470'
    "E44" = '  +4.35%  '
    "D45" = '2.96'
    "E45" = '  -0.37%  '
    "E46" = '  +4.55%  '
    "D47" = '3.39'
    "E47" = '  -2.64%  '
    "D48" = '9.05'
    "E48" = '  -4.48%  '
    "E49" = '  +0.52%  '
    "D50" = '133.99'
    "E50" = '  -0.82%  '
    "E51" = '  -5.36%  '
}

foreach ($cellRef in $updates.Keys) {
    Set-TextValue $cellRef $updates[$cellRef]
}

